$wb = $excel.ActiveWorkbook

# --- Update ENVELOPE_ASSEMBLIES data: add missing floor assemblies ---
$ws2 = $wb.Worksheets.Item("ENVELOPE_ASSEMBLIES")
# type_base column (I) - FLOOR_AS2 -> FLOOR_AS4 for all data rows
[void]($ws2.Range("I2:I7").Value = "FLOOR_AS4")
# type_floor column (H) - rows 3-7 FLOOR_AS1 -> FLOOR_AS3 (row 2 stays FLOOR_AS1)
[void]($ws2.Range("H3:H7").Value = "FLOOR_AS3")

# --- Update selection on STANDARD_DEFINITION sheet ---
$ws1 = $wb.Worksheets.Item("STANDARD_DEFINITION")
[void]$ws1.Activate()
[void]$ws1.Range("B3:B6").Select()

# --- Activate ENVELOPE_ASSEMBLIES sheet and update its selection ---
[void]$ws2.Activate()
[void]$ws2.Range("H3").Select()
